$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 26.144619
$ws.Range("H2").Value = 78.433857
$ws.Range("I2").Value = 0.5211737020083955
$ws.Range("J2").Value = 0.5211737020083955
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 13.10121233333333
$ws.Range("N2").Value = 39.303637
$ws.Range("O2").Value = 0.1081423012186565
$ws.Range("P2").Value = 0.1081423012186565
$ws.Range("Q2").Value = 342.5262048931011
$ws.Range("R2").Value = 3082.735844037909
$ws.Range("S2").Value = 0.05636092346983421
$ws.Range("T2").Value = 0.05636092346983422
$ws.Range("G3").Value = 26.144619
$ws.Range("H3").Value = 78.433857
$ws.Range("I3").Value = 0.5211737020083955
$ws.Range("J3").Value = 0.5211737020083955
$ws.Range("O3").Value = 0.5751439322003361
$ws.Range("P3").Value = 0.5751439322003362
$ws.Range("Q3").Value = 1821.691106476009
$ws.Range("R3").Value = 16395.21995828408
$ws.Range("S3").Value = 0.2997498923325148
$ws.Range("T3").Value = 0.2997498923325149
$ws.Range("G4").Value = 26.144619
$ws.Range("H4").Value = 78.433857
$ws.Range("I4").Value = 0.5211737020083955
$ws.Range("J4").Value = 0.5211737020083955
$ws.Range("M4").Value = 38.36920666666666
$ws.Range("N4").Value = 115.10762
$ws.Range("O4").Value = 0.3167137665810073
$ws.Range("P4").Value = 0.3167137665810074
$ws.Range("Q4").Value = 1003.14828963226
$ws.Range("R4").Value = 9028.33460669034
$ws.Range("S4").Value = 0.1650628862060464
$ws.Range("T4").Value = 0.1650628862060465
$ws.Range("I5").Value = 0.3571392594830743
$ws.Range("J5").Value = 0.3571392594830742
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 13.10121233333333
$ws.Range("N5").Value = 39.303637
$ws.Range("O5").Value = 0.1081423012186565
$ws.Range("P5").Value = 0.1081423012186565
$ws.Range("Q5").Value = 234.7193549821501
$ws.Range("R5").Value = 2112.47419483935
$ws.Range("S5").Value = 0.03862186137602653
$ws.Range("T5").Value = 0.03862186137602653
$ws.Range("I6").Value = 0.3571392594830743
$ws.Range("J6").Value = 0.3571392594830742
$ws.Range("O6").Value = 0.5751439322003361
$ws.Range("P6").Value = 0.5751439322003362
$ws.Range("S6").Value = 0.2054064780422115
$ws.Range("T6").Value = 0.2054064780422115
$ws.Range("I7").Value = 0.3571392594830743
$ws.Range("J7").Value = 0.3571392594830742
$ws.Range("M7").Value = 38.36920666666666
$ws.Range("N7").Value = 115.10762
$ws.Range("O7").Value = 0.3167137665810073
$ws.Range("P7").Value = 0.3167137665810074
$ws.Range("Q7").Value = 687.416951259
$ws.Range("R7").Value = 6186.752561331
$ws.Range("S7").Value = 0.1131109200648362
$ws.Range("T7").Value = 0.1131109200648362
$ws.Range("G8").Value = 6.104416333333333
$ws.Range("H8").Value = 18.313249
$ws.Range("I8").Value = 0.1216870385085301
$ws.Range("J8").Value = 0.1216870385085301
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 13.10121233333333
$ws.Range("N8").Value = 39.303637
$ws.Range("O8").Value = 0.1081423012186565
$ws.Range("P8").Value = 0.1081423012186565
$ws.Range("Q8").Value = 79.97525455406812
$ws.Range("R8").Value = 719.777290986613
$ws.Range("S8").Value = 0.01315951637279572
$ws.Range("T8").Value = 0.01315951637279572
$ws.Range("G9").Value = 6.104416333333333
$ws.Range("H9").Value = 18.313249
$ws.Range("I9").Value = 0.1216870385085301
$ws.Range("J9").Value = 0.1216870385085301
$ws.Range("O9").Value = 0.5751439322003361
$ws.Range("P9").Value = 0.5751439322003362
$ws.Range("Q9").Value = 425.3403327338685
$ws.Range("R9").Value = 3828.062994604816
$ws.Range("S9").Value = 0.06998756182560975
$ws.Range("T9").Value = 0.06998756182560976
$ws.Range("G10").Value = 6.104416333333333
$ws.Range("H10").Value = 18.313249
$ws.Range("I10").Value = 0.1216870385085301
$ws.Range("J10").Value = 0.1216870385085301
$ws.Range("M10").Value = 38.36920666666666
$ws.Range("N10").Value = 115.10762
$ws.Range("O10").Value = 0.3167137665810073
$ws.Range("P10").Value = 0.3167137665810074
$ws.Range("Q10").Value = 234.2216118730422
$ws.Range("R10").Value = 2107.99450685738
$ws.Range("S10").Value = 0.03853996031012467
$ws.Range("T10").Value = 0.03853996031012467
